$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain (unstyled) style reference used to avoid leaving a stray number-format
# style applied to cells that need to be forced to text so they are not
# reinterpreted as numbers (e.g. "1.00", "0.997", ...).
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "59.257.42"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.525.17"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.89"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.53"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "2.523.96"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").Value = "2.975.90"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.23"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "59.141.60"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.523.64"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.35"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.70"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +5.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.423"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.78"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.62"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +4.69%  "
$ws.Range("E33").Value = "  +6.80%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.54"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.68"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.825"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "284.39"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.27"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "130.73"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +7.52%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.88"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0932"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.54"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -0.14%  "
